$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.681.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.212.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.635'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.10'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.406'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0867'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.540.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.826'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.208.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '40.510.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '249.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.47%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -2.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '173.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("E34").Value = '  -1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.11'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("E37").Value = '  +6.01%  '
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.40%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.44%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.523.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0936'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.78%  '
$ws.Range("E48").Value = '  +40.45%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.70%  '
